$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells retain their original Text data type (they were inline strings),
# not auto-converted to numbers/percentages by Excel when re-entering numeric-looking text.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "297.63"
$ws.Range("E2").Value = "1.84%"
$ws.Range("D3").Value = "41.76"
$ws.Range("E3").Value = "3.48%"
$ws.Range("D4").Value = "5.013"
$ws.Range("E4").Value = "-0.13%"
$ws.Range("D5").Value = "0.07522"
$ws.Range("E5").Value = "3.15%"
$ws.Range("D6").Value = "1.601"
$ws.Range("E6").Value = "4.46%"
$ws.Range("D7").Value = "0.9160"
$ws.Range("E7").Value = "-1.08%"
$ws.Range("E8").Value = "2.18%"
$ws.Range("D9").Value = "0.1181"
$ws.Range("E9").Value = "1.65%"
$ws.Range("D10").Value = "0.1822"
$ws.Range("E10").Value = "3.40%"
$ws.Range("D11").Value = "0.08888"
$ws.Range("E11").Value = "2.15%"
$ws.Range("D12").Value = "0.04116"
$ws.Range("E12").Value = "-5.41%"
$ws.Range("D13").Value = "0.1051"
$ws.Range("E13").Value = "-0.14%"
$ws.Range("D14").Value = "0.001281"
$ws.Range("E14").Value = "0.31%"
$ws.Range("D15").Value = "0.006020"
$ws.Range("E15").Value = "1.00%"
$ws.Range("D16").Value = "3.342"
$ws.Range("E16").Value = "0.04%"
$ws.Range("D17").Value = "4.370"
$ws.Range("E17").Value = "1.91%"
$ws.Range("D18").Value = "0.3330"
$ws.Range("E18").Value = "1.48%"
$ws.Range("D19").Value = "8.325"
$ws.Range("E19").Value = "4.26%"
$ws.Range("E20").Value = "-2.90%"
$ws.Range("E21").Value = "11.78%"
$ws.Range("D22").Value = "0.04106"
$ws.Range("E22").Value = "4.69%"
$ws.Range("D23").Value = "0.001267"
$ws.Range("E23").Value = "0.27%"
$ws.Range("D24").Value = "0.003910"
$ws.Range("E24").Value = "6.77%"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").Value = "8.30%"
$ws.Range("E38").Value = "4.36%"
$ws.Range("D39").Value = "0.05217"
$ws.Range("E39").Value = "3.65%"
$ws.Range("D40").Value = "0.006312"
$ws.Range("E40").Value = "3.78%"
$ws.Range("D41").Value = "0.007776"
$ws.Range("E41").Value = "-1.03%"
$ws.Range("D42").Value = "0.1324"
$ws.Range("E42").Value = "3.14%"
$ws.Range("D43").Value = "0.007417"
$ws.Range("E43").Value = "0.77%"
$ws.Range("D44").Value = "0.007725"
$ws.Range("E44").Value = "6.48%"
$ws.Range("D45").Value = "0.3237"
$ws.Range("E45").Value = "1.70%"
$ws.Range("D46").Value = "0.00006589"
$ws.Range("E46").Value = "6.63%"
$ws.Range("E47").Value = "-0.02%"
$ws.Range("D48").Value = "0.04535"
$ws.Range("E48").Value = "-7.26%"
$ws.Range("D49").Value = "0.004208"
$ws.Range("E49").Value = "0.15%"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").Value = "-0.02%"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").Value = "-0.02%"
